# Committed Corporate Customer excel file
# Insert 4 new columns (CHEQUE.NUMBER, ORDERING.CUST:1, PAYMENT.DETAILS:1,
# COMMISSION.AMT:1) before the old "DD.ADDRESS:1" column, and append a new
# trailing "PURPOSE:1" column after the old last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift old columns D:H (DD.ADDRESS:1 .. INS.ISS.PURPOSE) right by four,
# opening up D:G for the new headers.
$ws.Range("D1:G1").EntireColumn.Insert()

# New columns get a plain 16-character custom width (no bestFit autofit).
$ws.Range("D1:G2").ColumnWidth = 15.1666666666667

# New header row values (row 2 under these columns is intentionally left
# blank, matching the source data).
$ws.Range("D1").Value = "CHEQUE.NUMBER"
$ws.Range("E1").Value = "ORDERING.CUST:1"
$ws.Range("F1").Value = "PAYMENT.DETAILS:1"
$ws.Range("G1").Value = "COMMISSION.AMT:1"

# New trailing header column (M), also with an empty data row underneath.
$ws.Range("M1").Value = "PURPOSE:1"

# Leave the selection where the author's last edit landed.
$ws.Range("M1").Select()
